# Improved logging and data caching. Using Adam optimizer now.
# Appends new experiment result rows (3-17) to the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 1).Value = "-x**2"
$ws.Cells.Item(3, 2).Value = 20000
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 5).Value = 0.08142967522144318
$ws.Cells.Item(3, 6).Value = 19.83478403091431
$ws.Cells.Item(3, 7).Value = 1008.329607664404

$ws.Cells.Item(4, 1).Value = "x**2"
$ws.Cells.Item(4, 2).Value = 20000
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0.000000000000002086819613799132
$ws.Cells.Item(4, 6).Value = 16.9700231552124
$ws.Cells.Item(4, 7).Value = 1178.548774923559

$ws.Cells.Item(5, 1).Value = "-x**2"
$ws.Cells.Item(5, 2).Value = 20000
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 0.00125144375488162
$ws.Cells.Item(5, 6).Value = 16.89823484420776
$ws.Cells.Item(5, 7).Value = 1183.555571596014

$ws.Cells.Item(6, 1).Value = "sin(10*x)"
$ws.Cells.Item(6, 2).Value = 20000
$ws.Cells.Item(6, 3).Value = 4
$ws.Cells.Item(6, 4).Value = 4
$ws.Cells.Item(6, 5).Value = 0.000135384892928414
$ws.Cells.Item(6, 6).Value = 79.72848176956177
$ws.Cells.Item(6, 7).Value = 250.8513840487487

$ws.Cells.Item(7, 1).Value = "-x**2"
$ws.Cells.Item(7, 2).Value = 20000
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 1
$ws.Cells.Item(7, 5).Value = 0.00000116468106625689
$ws.Cells.Item(7, 6).Value = 23.01174783706665
$ws.Cells.Item(7, 7).Value = 869.1212915076613

$ws.Cells.Item(8, 1).Value = "sin(5*x)"
$ws.Cells.Item(8, 2).Value = 20000
$ws.Cells.Item(8, 3).Value = 2
$ws.Cells.Item(8, 4).Value = 2
$ws.Cells.Item(8, 5).Value = 0.0003030607185792178
$ws.Cells.Item(8, 6).Value = 49.36098575592041
$ws.Cells.Item(8, 7).Value = 405.178294025483

$ws.Cells.Item(9, 1).Value = "sin(10*x)"
$ws.Cells.Item(9, 2).Value = 20000
$ws.Cells.Item(9, 3).Value = 4
$ws.Cells.Item(9, 4).Value = 4
$ws.Cells.Item(9, 5).Value = 0.000135384892928414
$ws.Cells.Item(9, 6).Value = 80.61733603477478
$ws.Cells.Item(9, 7).Value = 248.0855977599269

$ws.Cells.Item(10, 1).Value = "sin(10*x)"
$ws.Cells.Item(10, 2).Value = 20000
$ws.Cells.Item(10, 3).Value = 4
$ws.Cells.Item(10, 4).Value = 4
$ws.Cells.Item(10, 5).Value = 0.0002126327308360487
$ws.Cells.Item(10, 6).Value = 83.81179332733154
$ws.Cells.Item(10, 7).Value = 238.6299016641835

$ws.Cells.Item(11, 1).Value = "sin(10*x)"
$ws.Cells.Item(11, 2).Value = 20000
$ws.Cells.Item(11, 3).Value = 4
$ws.Cells.Item(11, 4).Value = 4
$ws.Cells.Item(11, 5).Value = 0.0002079298283206299
$ws.Cells.Item(11, 6).Value = 83.99549984931946
$ws.Cells.Item(11, 7).Value = 238.1079943077694

$ws.Cells.Item(12, 1).Value = "sin(10*x)"
$ws.Cells.Item(12, 2).Value = 20000
$ws.Cells.Item(12, 3).Value = 4
$ws.Cells.Item(12, 4).Value = 4
$ws.Cells.Item(12, 5).Value = 0.0002149988722521812
$ws.Cells.Item(12, 6).Value = 82.73853874206543
$ws.Cells.Item(12, 7).Value = 241.7253229761444

$ws.Cells.Item(13, 1).Value = "sin(10*x)"
$ws.Cells.Item(13, 2).Value = 20000
$ws.Cells.Item(13, 3).Value = 4
$ws.Cells.Item(13, 4).Value = 4
$ws.Cells.Item(13, 5).Value = 0.0002149988722521812
$ws.Cells.Item(13, 6).Value = 82.53280305862427
$ws.Cells.Item(13, 7).Value = 242.3278897457742

$ws.Cells.Item(14, 1).Value = "sin(10*x)"
$ws.Cells.Item(14, 2).Value = 20000
$ws.Cells.Item(14, 3).Value = 4
$ws.Cells.Item(14, 4).Value = 4
$ws.Cells.Item(14, 5).Value = 0.0002149988722521812
$ws.Cells.Item(14, 6).Value = 83.14774823188782
$ws.Cells.Item(14, 7).Value = 240.5356780585652

$ws.Cells.Item(15, 1).Value = "sin(10*x)"
$ws.Cells.Item(15, 2).Value = 20000
$ws.Cells.Item(15, 3).Value = 4
$ws.Cells.Item(15, 4).Value = 4
$ws.Cells.Item(15, 5).Value = 0.0002149988722521812
$ws.Cells.Item(15, 6).Value = 82.21243619918823
$ws.Cells.Item(15, 7).Value = 243.2721973053206

$ws.Cells.Item(16, 1).Value = "sin(10*x)"
$ws.Cells.Item(16, 2).Value = 20000
$ws.Cells.Item(16, 3).Value = 4
$ws.Cells.Item(16, 4).Value = 4
$ws.Cells.Item(16, 5).Value = 0.0002149988722521812
$ws.Cells.Item(16, 6).Value = 82.92414236068726
$ws.Cells.Item(16, 7).Value = 241.184285186911

$ws.Cells.Item(17, 1).Value = "sin(10*x)"
$ws.Cells.Item(17, 2).Value = 20000
$ws.Cells.Item(17, 3).Value = 4
$ws.Cells.Item(17, 4).Value = 4
$ws.Cells.Item(17, 5).Value = 0.0002149988722521812
$ws.Cells.Item(17, 6).Value = 80.24481797218323
$ws.Cells.Item(17, 7).Value = 249.2372779377851

